$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.078.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.441.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'580.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.37%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.15%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.439.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = "'  -4.00%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.32%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.38%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.344"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.99%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'26.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -4.78%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.875.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.05%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'62.168.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.430.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.00%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.17%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'329.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -6.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'65.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.35%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'618.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.78%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.562.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.65%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0943"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -9.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -6.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.95%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.91%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +0.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -6.86%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.374"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'149.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.73%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'42.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -9.28%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'142.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.71%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0523"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.91%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'19.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -8.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0₆0233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.63%  "
$ws.Range("E51").Style = "Normal"
